$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number line and the "Report Covering the Week"
# line). These shared strings are made of several formatting runs that all
# share identical formatting, so simply re-assigning the full text keeps the
# same look.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# ---------------------------------------------------------------------------
# Helper cells used as "format donors" so that when a cell needs to flip
# between a numeric value and the text placeholders ("0" / "***.*") it ends
# up with the same cell style the real workbook uses for that kind of cell.
#   C16 -> style used for text placeholder "0"
#   E18 -> style used for text placeholder "***.*"
#   C17 -> style used for plain whole numbers
# ---------------------------------------------------------------------------
$txtZero = $ws.Range("C16")
$txtDash = $ws.Range("E18")
$numStyle = $ws.Range("C17")
$pctStyle = $ws.Range("E16")

# ---------------------------------------------------------------------------
# Row 16 - Rape
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 1
$ws.Range("G16").Value = 4
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = -89.473684210526
$ws.Range("L16").Value = -50
$ws.Range("N16").Value = -95

# ---------------------------------------------------------------------------
# Row 17 - Robbery
# ---------------------------------------------------------------------------
$txtZero.Copy($ws.Range("D17"))
$txtDash.Copy($ws.Range("E17"))
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 6
$ws.Range("K17").Value = 50
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = -57.142857142857

# ---------------------------------------------------------------------------
# Row 19 - Burglary
# ---------------------------------------------------------------------------
$numStyle.Copy($ws.Range("C19"))
$ws.Range("C19").Value = 1
$txtZero.Copy($ws.Range("D19"))
$txtDash.Copy($ws.Range("E19"))
$ws.Range("F19").Value = 2
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 12
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = -14.285714285714
$ws.Range("N19").Value = -63.636363636363

# ---------------------------------------------------------------------------
# Row 21 - G.L.A. (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 100
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = -37.5
$ws.Range("I21").Value = 21
$ws.Range("J21").Value = 36
$ws.Range("K21").Value = -41.666666666666
$ws.Range("L21").Value = 31.25
$ws.Range("M21").Value = -8.695652173913
$ws.Range("N21").Value = -79.611650485436

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("F24").Value = 3
$ws.Range("H24").Value = 0
$ws.Range("M24").Value = -54.545454545454

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$txtZero.Copy($ws.Range("C26"))
$numStyle.Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$pctStyle.Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -42.857142857142

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$numStyle.Copy($ws.Range("C28"))
$ws.Range("C28").Value = 3
$numStyle.Copy($ws.Range("F28"))
$ws.Range("F28").Value = 3
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = -42.857142857142
